$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Grabs a paragraph's own "<w:p ...>" opening tag (with its w14:paraId /
# w:rsidR / etc. identity attributes intact) so we can reuse it verbatim
# when rebuilding that same paragraph via InsertXML, keeping the rest of
# the document's metadata untouched.
function Get-ParaOpenTag($paragraph) {
    $xml = $paragraph.Range.WordOpenXML
    if ($xml -match '<w:p( [^>]*)?>') {
        $attrs = $matches[1]
        if ($attrs -eq $null) { $attrs = "" }
        return "<w:p $wNs xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'$attrs>"
    }
    return "<w:p $wNs>"
}

# --- Step 1: insert a brand new "MainPicture NTEXT NOT NULL," paragraph
#     immediately before the existing "Title NVARCHAR(200) NOT NULL," line
#     in the News table definition. InsertParagraphBefore() first creates a
#     clean, empty paragraph without disturbing the Title paragraph; then
#     InsertXML fills it so the leading tab is a real <w:tab/> run child
#     (as Word represents it) instead of a literal tab character in <w:t>.
$rng = $d.Content
[void]$rng.Find.Execute("Title NVARCHAR(200) NOT NULL,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titlePara = $rng.Paragraphs(1)
$titleIndex = $titlePara.Index
[void]$titlePara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs($titleIndex)
[void]$newPara.Range.InsertXML("<w:p $wNs><w:r><w:tab/><w:t>MainPicture NTEXT NOT NULL,</w:t></w:r></w:p>")

# --- Step 2: the <w:lastRenderedPageBreak/> marker used to sit at the start
#     of the Description run; now that a paragraph was inserted above it, it
#     renders before the Title line instead, so move it onto Title's run
#     (preserving Title's own paragraph identity attributes).
$rng2 = $d.Content
[void]$rng2.Find.Execute("Title NVARCHAR(200) NOT NULL,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titlePara2 = $rng2.Paragraphs(1)
$titleTag = Get-ParaOpenTag $titlePara2
[void]$titlePara2.Range.InsertXML("$titleTag<w:r><w:lastRenderedPageBreak/><w:tab/><w:t>Title NVARCHAR(200) NOT NULL,</w:t></w:r></w:p>")

# --- Step 3: remove the now-stale <w:lastRenderedPageBreak/> from the
#     Description paragraph (preserving its own paragraph identity attrs).
$rng3 = $d.Content
[void]$rng3.Find.Execute("Description NTEXT NOT NULL,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$descPara = $rng3.Paragraphs(1)
$descTag = Get-ParaOpenTag $descPara
[void]$descPara.Range.InsertXML("$descTag<w:r><w:tab/><w:t>Description NTEXT NOT NULL,</w:t></w:r></w:p>")

Write-Host "MainPicture column inserted before Title; lastRenderedPageBreak moved."
